$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# These values correspond to a reshuffle of the weekly data rows.
$rowData = @{
    2  = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
    3  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    4  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    5  = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    6  = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
    7  = @{ D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 }
    8  = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    9  = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    10 = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
    11 = @{ D = 45068; J = 400; K = 16000; L = 17000; M = 16500; P = 917 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
